$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "68.087.68"
$ws.Range("E2").Value = "  -0.72%  "
$ws.Range("D3").Value = "3.869.51"
$ws.Range("E3").Value = "  -1.10%  "
$ws.Range("D4").Value = "'1.00"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  -0.08%  "
$ws.Range("D5").Value = "'599.43"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -0.70%  "
$ws.Range("D6").Value = "'167.34"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +1.01%  "
$ws.Range("D7").Value = "3.868.04"
$ws.Range("E7").Value = "  -1.10%  "
$ws.Range("E8").Value = "  +0.06%  "
$ws.Range("E9").Value = "  -0.71%  "
$ws.Range("E10").Value = "  -0.54%  "
$ws.Range("D11").Value = "'6.41"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +1.10%  "
$ws.Range("D12").Value = "'0.457"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -0.60%  "
$ws.Range("E13").Value = "  -0.41%  "
$ws.Range("D14").Value = "'36.94"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -0.79%  "
$ws.Range("D15").Value = "4.527.12"
$ws.Range("E15").Value = "  -0.90%  "
$ws.Range("D16").Value = "3.873.04"
$ws.Range("E16").Value = "  -1.41%  "
$ws.Range("D17").Value = "68.130.30"
$ws.Range("E17").Value = "  -0.86%  "
$ws.Range("D18").Value = "'18.26"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +6.43%  "
$ws.Range("D19").Value = "'7.41"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -0.54%  "
$ws.Range("E20").Value = "  -1.12%  "
$ws.Range("D21").Value = "'10.89"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -1.16%  "
$ws.Range("D22").Value = "'465.83"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -4.09%  "
$ws.Range("D23").Value = "'0.732"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +1.11%  "
$ws.Range("E24").Value = "  -3.34%  "
$ws.Range("D25").Value = "'83.50"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -1.48%  "
$ws.Range("E26").Value = "  +0.57%  "
$ws.Range("D27").Value = "'12.12"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +0.75%  "
$ws.Range("D28").Value = "'10.01"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -1.17%  "
$ws.Range("E29").Value = "  +0.05%  "
$ws.Range("D30").Value = "'2.95"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +0.65%  "
$ws.Range("D31").Value = "4.019.50"
$ws.Range("E31").Value = "  -1.22%  "
$ws.Range("D32").Value = "'7.85"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +0.60%  "
$ws.Range("E33").Value = "  -3.36%  "
$ws.Range("D34").Value = "'31.24"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -2.45%  "
$ws.Range("D35").Value = "'9.35"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +1.21%  "
$ws.Range("D36").Value = "3.844.94"
$ws.Range("E36").Value = "  -0.45%  "
$ws.Range("E37").Value = "  -2.68%  "
$ws.Range("D38").Value = "'3.43"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +8.58%  "
$ws.Range("E39").Value = "  -1.51%  "
$ws.Range("E40").Value = "  +0.88%  "
$ws.Range("E41").Value = "  +0.15%  "
$ws.Range("D42").Value = "'1.00"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +0.08%  "
$ws.Range("D43").Value = "'0.313"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -1.38%  "
$ws.Range("D44").Value = "'429.05"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +0.23%  "
$ws.Range("E45").Value = "  -0.19%  "
$ws.Range("E46").Value = "  -0.01%  "
$ws.Range("D47").Value = "'47.33"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -2.35%  "
$ws.Range("D48").Value = "'8.53"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +0.34%  "
$ws.Range("D49").Value = "'0.000278"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +5.76%  "
$ws.Range("B50").Value = "Arweave"
$ws.Range("C50").Value = "https://coinranking.com/coin/7XWg41D1+arweave-ar"
$ws.Range("D50").Value = "'40.71"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +3.66%  "
$ws.Range("B51").Value = "Monero"
$ws.Range("C51").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D51").Value = "'144.23"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +1.44%  "
